# no-op
$wb = $excel.ActiveWorkbook
